$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns before column C (old C shifts to E; new C & D are
#    inserted blank in between). This mirrors "add this week's rating column"
#    for the MarketBeat tracker: B keeps the newest week, C/D are the two new
#    weeks being backfilled with "UN" (unrated), and the old C (now E) keeps
#    its original per-firm rating text/highlight.
# ---------------------------------------------------------------------------
$ws.Columns("C:D").Insert()

# Match column widths (same custom width as the original column C == 8.0)
$ws.Range("C1").ColumnWidth = 7.166666666666667
$ws.Range("D1").ColumnWidth = 7.166666666666667
$ws.Range("E1").ColumnWidth = 7.166666666666667

# ---------------------------------------------------------------------------
# 2. Header row: shift the date labels right and add the two new week labels.
#    Old C1 ("Jun_10") already moved into E1 via the column insert/shift.
#    B1 wasn't touched by the insert (it's left of the insertion point), so
#    its old value ("Jun_13") is carried over into D1 by hand, then B1/C1
#    get the two new week labels.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Jun_13"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# ---------------------------------------------------------------------------
# 3. Fill the two new columns (C, D) with "UN" (unrated) for every data row.
# ---------------------------------------------------------------------------
$ws.Range("C2:D27").Value = "UN"

# ---------------------------------------------------------------------------
# 4. Row 22 (BidaskClub) got an actual new rating headline in the newest
#    column (B) instead of "UN", highlighted the same way the other
#    standout cells (e.g. old C18, now E18) are.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "6/16/2018,Upgrades,Hold -> Buy,"
$ws.Range("B22").Interior.ColorIndex = 35

Write-Output "done"
